$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Φύλλο1")

# Copy the formatting of row 15 down onto row 16 so new cells pick up the
# same fonts/alignment used by the rest of the transactions table, then
# fill in the new "giveRightToVote" transaction for Eniola.
$ws.Range("A15:F15").Copy()
$ws.Range("A16:F16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A16").Value = "0x325fde66e12f6b44a18a1512e9cc138f5fa705a931fb6b6274d7ff8630ac3202"
$ws.Range("B16").Value = "0x56814ca0854e878c2fd9ffa0899c12f4c4e35346"
$ws.Range("C16").Value = "0x41874b62017e8bf2c533b73c11750fdbb78ac956"
$ws.Range("D16").Value = "giveRightToVote / address / `t0xE3A9a11232f4D52786CA61f56bB7Fb01b00C80cd"
$ws.Range("E16").Value = "Success"
$ws.Range("F16").Value = "Nikos give rights to to Eniola so he can vote"

# Restore the previous selection/scroll state that Excel records on save.
$ws.Range("B6").Select()
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Application.ActiveWindow.ScrollRow = 1
